$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "59.862.31"
$c.Style = $origStyle
$ws.Cells.Item(2, 5).Value = "  -1.86%  "

# Row 3
$c = $ws.Cells.Item(3, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.373.17"
$c.Style = $origStyle
$ws.Cells.Item(3, 5).Value = "  -2.47%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.10%  "

# Row 5
$c = $ws.Cells.Item(5, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "559.95"
$c.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  -2.06%  "

# Row 6
$c = $ws.Cells.Item(6, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "138.37"
$c.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  -1.64%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.13%  "

# Row 8
$c = $ws.Cells.Item(8, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.527"
$c.Style = $origStyle
$ws.Cells.Item(8, 5).Value = "  -0.48%  "

# Row 9
$c = $ws.Cells.Item(9, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.369.44"
$c.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  -2.10%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -3.93%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -1.17%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -1.12%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.84%  "

# Row 14
$c = $ws.Cells.Item(14, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "25.53"
$c.Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  -2.21%  "

# Row 15
$c = $ws.Cells.Item(15, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.796.72"
$c.Style = $origStyle
$ws.Cells.Item(15, 5).Value = "  -1.12%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -3.34%  "

# Row 17
$c = $ws.Cells.Item(17, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "59.633.40"
$c.Style = $origStyle
$ws.Cells.Item(17, 5).Value = "  -2.11%  "

# Row 18
$c = $ws.Cells.Item(18, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.382.31"
$c.Style = $origStyle
$ws.Cells.Item(18, 5).Value = "  -1.28%  "

# Row 19
$c = $ws.Cells.Item(19, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.09"
$c.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  +10.95%  "

# Row 20
$c = $ws.Cells.Item(20, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.45"
$c.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  -1.23%  "

# Row 21
$c = $ws.Cells.Item(21, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "321.52"
$c.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  -0.69%  "

# Row 22
$c = $ws.Cells.Item(22, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  -0.22%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -2.84%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  +0.02%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -4.22%  "

# Row 26
$c = $ws.Cells.Item(26, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "64.15"
$c.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  -1.05%  "

# Row 27
$c = $ws.Cells.Item(27, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "561.45"
$c.Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  -2.77%  "

# Row 28
$c = $ws.Cells.Item(28, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.12"
$c.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  -8.08%  "

# Row 29
$c = $ws.Cells.Item(29, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.483.39"
$c.Style = $origStyle
$ws.Cells.Item(29, 5).Value = "  -2.99%  "

# Row 30
$c = $ws.Cells.Item(30, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0₃0923"
$c.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  +0.97%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +1.29%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -3.62%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -3.45%  "

# Row 34
$c = $ws.Cells.Item(34, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.130"
$c.Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  -1.90%  "

# Row 35
$c = $ws.Cells.Item(35, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = $origStyle
$ws.Cells.Item(35, 5).Value = "  -0.49%  "

# Row 36
$c = $ws.Cells.Item(36, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  +2.50%  "

# Row 37
$c = $ws.Cells.Item(37, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "152.69"
$c.Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  +1.66%  "

# Row 38
$c = $ws.Cells.Item(38, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.366"
$c.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  -0.34%  "

# Row 39
$c = $ws.Cells.Item(39, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "4.52"
$c.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  -1.73%  "

# Row 40
$c = $ws.Cells.Item(40, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.10"
$c.Style = $origStyle
$ws.Cells.Item(40, 5).Value = "  -0.89%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -2.41%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  -0.03%  "

# Row 43
$c = $ws.Cells.Item(43, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "41.52"
$c.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  -0.44%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.05%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "dogwifhat"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(45, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.39"
$c.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +1.92%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(46, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0₆0299"
$c.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  +6.54%  "

# Row 47
$c = $ws.Cells.Item(47, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "138.74"
$c.Style = $origStyle
$ws.Cells.Item(47, 5).Value = "  -1.77%  "

# Row 48
$c = $ws.Cells.Item(48, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.51"
$c.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  +0.01%  "

# Row 49
$c = $ws.Cells.Item(49, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.585"
$c.Style = $origStyle
$ws.Cells.Item(49, 5).Value = "  -1.51%  "

# Row 50
$c = $ws.Cells.Item(50, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0500"
$c.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  -1.28%  "

# Row 51
$c = $ws.Cells.Item(51, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "19.11"
$c.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  -2.25%  "
